# Fix a few issues with the ratios dataframe:
#  - insert a new "Lab. #" column at the very left (column A)
#  - fill in the Lab # values for every data row
#  - highlight the rows that belong to Lab # 10815 with a light green fill

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a brand-new column before the current column A. This shifts all the
# existing data (and column widths) one column to the right automatically.
$ws.Columns.Item(1).Insert()

# New header cell for the inserted column.
$ws.Cells.Item(1, 1).Value = "Lab. #"

# New column A is narrower than the data columns (target stored width is
# 7.7109375 characters; feed the engine's column-width rounding the value
# that snaps to the closest achievable width).
$ws.Columns.Item(1).ColumnWidth = 6.877604166666667

# Lab # values for each data row (rows 2-15).
$labNumbers = @{
    2  = 10815
    3  = 10973
    4  = 10815
    5  = 10974
    6  = 10815
    7  = 10975
    8  = 10815
    9  = 10976
    10 = 10815
    11 = 10977
    12 = 10815
    13 = 10978
    14 = 10815
    15 = 10979
}

foreach ($row in $labNumbers.Keys) {
    $ws.Cells.Item($row, 1).Value = $labNumbers[$row]
}

# Highlight every row whose Lab # is 10815 with a light green fill across
# all 9 columns (A through I).
$highlightRows = @(2, 4, 6, 8, 10, 12, 14)
foreach ($row in $highlightRows) {
    $rowRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 9))
    $rowRange.Interior.Color = 12379352
}
